$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.082.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.14%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.761.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.66%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.05%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''576.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.86%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''159.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.29%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.09%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -3.15%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -3.38%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = '''Toncoin'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '''5.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -13.79%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = '''TRON'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = '''https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '''0.165'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +3.34%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.386'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '''3.249.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.75%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''27.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.36%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''63.668.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.64%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  -5.49%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.765.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.78%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''12.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = '''4.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -3.92%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''359.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.36%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''6.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -6.13%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  -0.37%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.530'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -8.30%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''65.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -3.67%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -3.70%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''8.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -2.95%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.13%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.0₃0906'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -6.58%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''7.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +1.13%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -4.17%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''1.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +2.68%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''170.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -1.23%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''20.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -3.09%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''4.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -3.70%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = '''ImmutableX'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = '''1.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -1.22%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = '''USDe'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.09%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.85%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.59%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''350.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.18%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''6.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.98%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''4.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -2.82%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''39.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -2.23%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''21.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.31%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''21.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -3.79%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -4.39%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''137.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.21%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''VeChain'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''0.0254'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -3.17%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = '''Mantle'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''0.630'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -3.72%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -2.07%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  +0.01%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''11.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.11%  '
$ws.Range("E51").Style = "Normal"
